$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 42607.887233796297
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 52
$ws.Range("D5").Value = 47
$ws.Range("E5").Value = 66
$ws.Range("F5").Value = 33
$ws.Range("G5").Value = 10751
$ws.Range("H5").Value = 21557
$ws.Range("I5").Value = 2335
$ws.Range("J5").Value = 275
$ws.Range("K5").Value = 248
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = "Bag"

$ws.Range("A5").NumberFormat = "m/d/yy h:mm"
